$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

# Row -> (B value, C value)
$data = @{
    2  = @(2496.5, 2494.05)
    3  = @(382.25, 383.6)
    4  = @(1518.2, 1510.45)
    5  = @(7387.8, 7345.1)
    6  = @(234, 234.65)
    7  = @(195, 195.8)
    8  = @(44644.05, 44766.05)
    9  = @(519, 521.3)
    10 = @(3391.3, 3403.45)
    11 = @(143.8, 144.7)
    12 = @(1235.85, 1235.55)
    13 = @(1408.65, 1423.75)
    14 = @(702.45, 712.55)
    15 = @(449.25, 451.05)
    16 = @(1577, 1571.15)
    17 = @(294.35, 292.65)
    18 = @(19658.85, 19653.55)
    19 = @(570.5, 573)
    20 = @(615.35, 615.25)
    21 = @(610.25, 611.15)
    22 = @(257.5, 258.45)
    23 = @(129.6, 131.7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
